$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '59.613.11'
$ws.Cells.Item(2, 5).Value = '  +1.46%  '
$ws.Cells.Item(3, 4).Value = '2.587.34'
$ws.Cells.Item(3, 5).Value = '  +0.63%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
Set-TextValue 5 4 '555.99'
$ws.Cells.Item(5, 5).Value = '  -0.97%  '
Set-TextValue 6 4 '141.17'
$ws.Cells.Item(6, 5).Value = '  -1.25%  '
$ws.Cells.Item(7, 5).Value = '  -0.13%  '
$ws.Cells.Item(8, 5).Value = '  -0.29%  '
$ws.Cells.Item(9, 4).Value = '2.603.96'
$ws.Cells.Item(9, 5).Value = '  +1.09%  '
Set-TextValue 10 4 '6.69'
$ws.Cells.Item(10, 5).Value = '  +0.33%  '
$ws.Cells.Item(11, 5).Value = '  +1.23%  '
$ws.Cells.Item(12, 5).Value = '  +6.68%  '
Set-TextValue 13 4 '0.366'
$ws.Cells.Item(13, 5).Value = '  +7.52%  '
$ws.Cells.Item(14, 4).Value = '3.042.06'
$ws.Cells.Item(14, 5).Value = '  +0.69%  '
$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 15 4 '23.33'
$ws.Cells.Item(15, 5).Value = '  +5.85%  '
$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value = '59.590.19'
$ws.Cells.Item(16, 5).Value = '  +1.24%  '
Set-TextValue 17 4 '0.0000137'
$ws.Cells.Item(17, 5).Value = '  +0.70%  '
$ws.Cells.Item(18, 4).Value = '2.596.36'
$ws.Cells.Item(18, 5).Value = '  +0.84%  '
$ws.Cells.Item(19, 5).Value = '  +2.00%  '
Set-TextValue 20 4 '340.71'
$ws.Cells.Item(20, 5).Value = '  +1.68%  '
Set-TextValue 21 4 '10.52'
$ws.Cells.Item(21, 5).Value = '  +3.57%  '
Set-TextValue 22 4 '6.70'
$ws.Cells.Item(22, 5).Value = '  +9.22%  '
Set-TextValue 23 4 '0.996'
$ws.Cells.Item(23, 5).Value = '  -0.27%  '
Set-TextValue 24 4 '0.503'
$ws.Cells.Item(24, 5).Value = '  +11.52%  '
Set-TextValue 25 4 '62.39'
$ws.Cells.Item(25, 5).Value = '  -2.01%  '
Set-TextValue 26 4 '1.00'
$ws.Cells.Item(26, 5).Value = '  +0.02%  '
Set-TextValue 27 4 '0.159'
$ws.Cells.Item(27, 5).Value = '  -0.91%  '
$ws.Cells.Item(28, 5).Value = '  +3.91%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0774'
$ws.Cells.Item(29, 5).Value = '  -0.49%  '
$ws.Cells.Item(30, 5).Value = '  -0.11%  '
Set-TextValue 31 4 '1.69'
$ws.Cells.Item(31, 5).Value = '  +1.57%  '
$ws.Cells.Item(32, 5).Value = '  +1.92%  '
Set-TextValue 33 4 '158.55'
Set-TextValue 34 4 '19.30'
$ws.Cells.Item(34, 5).Value = '  +1.73%  '
Set-TextValue 35 4 '4.09'
$ws.Cells.Item(35, 5).Value = '  +2.45%  '
Set-TextValue 36 4 '0.914'
$ws.Cells.Item(36, 5).Value = '  +4.10%  '
$ws.Cells.Item(37, 5).Value = '  +4.27%  '
Set-TextValue 38 4 '37.74'
$ws.Cells.Item(38, 5).Value = '  +2.80%  '
$ws.Cells.Item(39, 5).Value = '  +2.11%  '
Set-TextValue 40 4 '0.842'
$ws.Cells.Item(40, 5).Value = '  -3.91%  '
$ws.Cells.Item(41, 5).Value = '  +1.50%  '
Set-TextValue 42 4 '290.71'
$ws.Cells.Item(42, 5).Value = '  +0.36%  '
Set-TextValue 43 4 '135.80'
$ws.Cells.Item(43, 5).Value = '  +10.06%  '
Set-TextValue 44 4 '0.997'
$ws.Cells.Item(44, 5).Value = '  -0.27%  '
Set-TextValue 45 4 '0.0976'
$ws.Cells.Item(45, 5).Value = '  +0.57%  '
$ws.Cells.Item(46, 5).Value = '  +1.16%  '
Set-TextValue 47 4 '0.0238'
$ws.Cells.Item(47, 5).Value = '  +3.33%  '
$ws.Cells.Item(48, 5).Value = '  +1.05%  '
Set-TextValue 49 4 '10.63'
$ws.Cells.Item(49, 5).Value = '  +0.12%  '
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 50 4 '18.85'
$ws.Cells.Item(50, 5).Value = '  +1.90%  '
$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).Value = '1.972.75'
$ws.Cells.Item(51, 5).Value = '  +2.38%  '
